$d = $word.ActiveDocument

# 1) Remove the trailing space after "Email," in the NHANVIEN(...) paragraph.
$null = $d.Content.Find.Execute("GioiTinh, SoDT, Email, ", $true, $false, $false, $false, $false, $true, 1, $false, "GioiTinh, SoDT, Email,", 2)

# 2) Insert " Username, Password," right after "Email," (new run, same formatting).
$rng1 = $d.Content
$null = $rng1.Find.Execute("Email,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Collapse(0)
$rng1.InsertAfter(" Username, Password,")
$afterUsernamePos = $rng1.End

# 3) Move the _GoBack bookmark from the end of the SANPHAM(...) paragraph to
#    right after "Username, Password," in the NHANVIEN(...) paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$bmRng = $d.Range($afterUsernamePos, $afterUsernamePos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# 4) Insert a single space before "MaLoaiNV" (after the bookmark).
$rng2 = $d.Range($afterUsernamePos, $afterUsernamePos)
$rng2.InsertAfter(" ")
